$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08)
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45207
}
